$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 23812360
$ws.Range("I6").Value = 55555730
$ws.Range("K6").Value = 166667190
$ws.Range("M6").Value = -166667078
$ws.Range("H43").Value = 9996.143
$ws.Range("I43").Value = 2000
$ws.Range("J43").Value = 13194.6
$ws.Range("K43").Value = 2000
$ws.Range("L43").Value = 13194.6
$ws.Range("M43").Value = -1931
$ws.Range("N43").Value = -13332.6
$ws.Range("H116").Value = 4498
$ws.Range("I116").Value = 5000
$ws.Range("J116").Value = 3996
$ws.Range("K116").Value = 5000
$ws.Range("L116").Value = 3996
$ws.Range("M116").Value = -1558
$ws.Range("N116").Value = -10880
$ws.Range("H129").Value = 2067.6667
$ws.Range("I129").Value = 990.6667
$ws.Range("J129").Value = 2426.6667
$ws.Range("K129").Value = 2972.0001
$ws.Range("L129").Value = 7280.000100000001
$ws.Range("M129").Value = 2027.9999
$ws.Range("N129").Value = -17280.0001
$ws.Range("H132").Value = 1615.4166
$ws.Range("I132").Value = 1465.2188
$ws.Range("J132").Value = 2817
$ws.Range("K132").Value = 4395.6564
$ws.Range("L132").Value = 8451
$ws.Range("M132").Value = -1865.6564
$ws.Range("N132").Value = -13511
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4316.6
$ws.Range("I2").Value = 4236.5454
$ws.Range("J2").Value = 4536.75
$ws.Range("K2").Value = 4236.5454
$ws.Range("L2").Value = 4536.75
$ws.Range("M2").Value = -4123.5454
$ws.Range("N2").Value = -4762.75
$ws.Range("H32").Value = 5424.92
$ws.Range("I32").Value = 5649.143
$ws.Range("K32").Value = 5649.143
$ws.Range("M32").Value = -5362.143
$ws.Range("H61").Value = 5185
$ws.Range("I61").Value = 5216.4443
$ws.Range("J61").Value = 4996.3335
$ws.Range("K61").Value = 5216.4443
$ws.Range("L61").Value = 4996.3335
$ws.Range("M61").Value = -5004.4443
$ws.Range("N61").Value = -5420.3335
$ws.Range("H74").Value = 2824.9285
$ws.Range("I74").Value = 2811.4614
$ws.Range("K74").Value = 2811.4614
$ws.Range("M74").Value = -1937.4614
$ws.Range("H77").Value = 2824.9285
$ws.Range("I77").Value = 2811.4614
$ws.Range("K77").Value = 14057.307
$ws.Range("M77").Value = -9689.307000000001
$ws.Range("H102").Value = 2657.9473
$ws.Range("I102").Value = 2250.0557
$ws.Range("K102").Value = 2250.0557
$ws.Range("M102").Value = -628.0556999999999
$ws.Range("H116").Value = 4316.6
$ws.Range("I116").Value = 4236.5454
$ws.Range("J116").Value = 4536.75
$ws.Range("K116").Value = 4236.5454
$ws.Range("L116").Value = 4536.75
$ws.Range("M116").Value = -1942.5454
$ws.Range("N116").Value = -9124.75
$ws.Range("H132").Value = 2548.5
$ws.Range("I132").Value = 2036
$ws.Range("K132").Value = 6108
$ws.Range("M132").Value = -3578
$ws.Range("H136").Value = 5185
$ws.Range("I136").Value = 5216.4443
$ws.Range("J136").Value = 4996.3335
$ws.Range("K136").Value = 15649.3329
$ws.Range("L136").Value = 14989.0005
$ws.Range("M136").Value = -13099.3329
$ws.Range("N136").Value = -20089.0005
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4316.6
$ws.Range("I3").Value = 4236.5454
$ws.Range("J3").Value = 4536.75
$ws.Range("K3").Value = 4236.5454
$ws.Range("L3").Value = 4536.75
$ws.Range("M3").Value = -4122.5454
$ws.Range("N3").Value = -4764.75
$ws.Range("H94").Value = 4286.636
$ws.Range("I94").Value = 3715.3
$ws.Range("K94").Value = 3715.3
$ws.Range("M94").Value = -3264.3
$ws.Range("H105").Value = 2579.4119
$ws.Range("I105").Value = 2450
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 2450
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -703
$ws.Range("N105").Value = -6494
$ws.Range("H107").Value = 5750
$ws.Range("I107").Value = 1500
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = 420
$ws.Range("N107").Value = -13840
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 170
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2346.6924
$ws.Range("I58").Value = 2333.9167
$ws.Range("K58").Value = 2333.9167
$ws.Range("M58").Value = -2130.9167
$ws.Range("H93").Value = 17720.111
$ws.Range("I93").Value = 13878.875
$ws.Range("K93").Value = 13878.875
$ws.Range("M93").Value = -12006.875
$ws.Range("H107").Value = 1652.2593
$ws.Range("I107").Value = 855.5
$ws.Range("J107").Value = 1716
$ws.Range("K107").Value = 855.5
$ws.Range("L107").Value = 1716
$ws.Range("M107").Value = 1064.5
$ws.Range("N107").Value = -5556
$ws.Range("H132").Value = 3002
$ws.Range("I132").Value = 2835.7222
$ws.Range("J132").Value = 3999.6667
$ws.Range("K132").Value = 8507.1666
$ws.Range("L132").Value = 11999.0001
$ws.Range("M132").Value = -5977.1666
$ws.Range("N132").Value = -17059.0001
$ws.Range("H134").Value = 3016.8235
$ws.Range("I134").Value = 2955.1428
$ws.Range("J134").Value = 3304.6667
$ws.Range("K134").Value = 8865.428400000001
$ws.Range("L134").Value = 9914.000100000001
$ws.Range("M134").Value = -6330.428400000001
$ws.Range("N134").Value = -14984.0001
$ws.Range("H136").Value = 2346.6924
$ws.Range("I136").Value = 2333.9167
$ws.Range("K136").Value = 7001.750100000001
$ws.Range("M136").Value = -4451.750100000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 1182.25
$ws.Range("J13").Value = 3721
$ws.Range("L13").Value = 11163
$ws.Range("N13").Value = -11499
$ws.Range("H82").Value = 47933.332
$ws.Range("I82").Value = 19000
$ws.Range("J82").Value = 62400
$ws.Range("K82").Value = 57000
$ws.Range("L82").Value = 187200
$ws.Range("M82").Value = -56594
$ws.Range("N82").Value = -188012
$ws.Range("H85").Value = 47933.332
$ws.Range("I85").Value = 19000
$ws.Range("J85").Value = 62400
$ws.Range("K85").Value = 57000
$ws.Range("L85").Value = 187200
$ws.Range("M85").Value = -55596
$ws.Range("N85").Value = -190008
$ws.Range("H92").Value = 3166
$ws.Range("J92").Value = 3498.6667
$ws.Range("L92").Value = 10496.0001
$ws.Range("N92").Value = -12992.0001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3693.2144
$ws.Range("I80").Value = 1978.3334
$ws.Range("J80").Value = 6780
$ws.Range("K80").Value = 1978.3334
$ws.Range("L80").Value = 6780
$ws.Range("M80").Value = -980.3334
$ws.Range("N80").Value = -8776
$ws.Range("H83").Value = 3693.2144
$ws.Range("I83").Value = 1978.3334
$ws.Range("J83").Value = 6780
$ws.Range("K83").Value = 9891.666999999999
$ws.Range("L83").Value = 33900
$ws.Range("M83").Value = -4899.666999999999
$ws.Range("N83").Value = -43884
$ws.Range("H102").Value = 2108.7334
$ws.Range("I102").Value = 2250.6667
$ws.Range("J102").Value = 1777.5555
$ws.Range("K102").Value = 2250.6667
$ws.Range("L102").Value = 1777.5555
$ws.Range("M102").Value = -628.6667000000002
$ws.Range("N102").Value = -5021.5555
$ws.Range("H113").Value = 8409.471
$ws.Range("I113").Value = 2602.2
$ws.Range("J113").Value = 10829.167
$ws.Range("K113").Value = 2602.2
$ws.Range("L113").Value = 10829.167
$ws.Range("M113").Value = -432.1999999999998
$ws.Range("N113").Value = -15169.167
$ws.Range("H139").Value = 99488
$ws.Range("J139").Value = 99488
$ws.Range("L139").Value = 99488
$ws.Range("N139").Value = -109768
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3253.2307
$ws.Range("I22").Value = 3194.6
$ws.Range("J22").Value = 3289.875
$ws.Range("K22").Value = 3194.6
$ws.Range("L22").Value = 3289.875
$ws.Range("M22").Value = -2899.6
$ws.Range("N22").Value = -3879.875
$ws.Range("H27").Value = 3253.2307
$ws.Range("I27").Value = 3194.6
$ws.Range("J27").Value = 3289.875
$ws.Range("K27").Value = 3194.6
$ws.Range("L27").Value = 3289.875
$ws.Range("M27").Value = -3087.6
$ws.Range("N27").Value = -3503.875
$ws.Range("H82").Value = 265
$ws.Range("I82").Value = 265
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 265
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = 96
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 265
$ws.Range("I85").Value = 265
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 265
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 983
$ws.Range("N85").ClearContents()
$ws.Range("H132").Value = 5564.913
$ws.Range("I132").Value = 5393.778
$ws.Range("J132").Value = 6181
$ws.Range("K132").Value = 16181.334
$ws.Range("L132").Value = 18543
$ws.Range("M132").Value = -13651.334
$ws.Range("N132").Value = -23603
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 39999.5
$ws.Range("J70").Value = 39999.5
$ws.Range("L70").Value = 39999.5
$ws.Range("N70").Value = -40629.5
$ws.Range("H73").Value = 39999.5
$ws.Range("J73").Value = 39999.5
$ws.Range("L73").Value = 39999.5
$ws.Range("N73").Value = -42183.5
$ws.Range("H122").Value = 4767.304
$ws.Range("I122").Value = 1356.125
$ws.Range("J122").Value = 6586.6
$ws.Range("K122").Value = 4068.375
$ws.Range("L122").Value = 19759.8
$ws.Range("M122").Value = -1618.375
$ws.Range("N122").Value = -24659.8
